$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, pushing the existing rows 39-43 down to 40-44
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly price-report record
$ws.Cells.Item(39, 1).Value  = 8
$ws.Cells.Item(39, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(39, 3).Value  = "Coquimbo"
$ws.Cells.Item(39, 4).Value  = 44476
$ws.Cells.Item(39, 5).Value  = 4
$ws.Cells.Item(39, 6).Value  = 100112052
$ws.Cells.Item(39, 7).Value  = "Albahaca"
$ws.Cells.Item(39, 8).Value  = "Sin especificar"
$ws.Cells.Item(39, 9).Value  = "Primera"
$ws.Cells.Item(39, 10).Value = 600
$ws.Cells.Item(39, 11).Value = 3500
$ws.Cells.Item(39, 12).Value = 4000
$ws.Cells.Item(39, 13).Value = 3750
$ws.Cells.Item(39, 14).Value = "`$/paquete"
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value = 3750
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
